# Israel Premier League - base update (27-03-2024 20:23)
#
# The source feed re-sorted same-kickoff-time fixtures and refreshed their
# results/odds, and appended the fixtures/odds for several newly scheduled
# matches. This script re-applies that refreshed snapshot on top of the
# existing workbook:
#   - re-writes the rows for kickoffs that share a timestamp with another
#     fixture, since the two rows traded places in the refreshed export
#   - fills in the final score/odds for two fixtures that had since been
#     played (rows 182 and 183)
#   - appends five brand-new fixtures at the end of the table (rows 184-188)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [object[]]$Values
    )
    $lastCol = $Values.Length
    $arr = New-Object 'object[,]' 1, $lastCol
    for ($i = 0; $i -lt $lastCol; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $startCell = $ws.Cells.Item($Row, 1)
    $endCell = $ws.Cells.Item($Row, $lastCol)
    $ws.Range($startCell, $endCell).Value = $arr
}

# --- Fixtures that swapped order with their same-kickoff-time sibling ---
Set-RowValues 4 @(2, 6798412, "Israel Premier League", "Israel Premier League", 45164.58333333334, "Hapoel Beer Sheva", "Hapoel Hadera", 3, 0, "H", 1.3, 4.5, 8, 1.333, 4.333, 7, -1.25, 1.825, 2.025, 2.75, 1.95, 1.9, 0.333, -1, -1, 0.825, -1, 0.475, -0.5)
Set-RowValues 5 @(3, 6799822, "Israel Premier League", "Israel Premier League", 45164.58333333334, "Maccabi Netanya", "Maccabi Bnei Raina", 1, 1, "D", 1.85, 3.3, 3.7, 1.65, 3.5, 4.5, -0.75, 1.95, 1.9, 2.5, 2, 1.85, -1, 2.5, -1, -1, 0.8999999999999999, -1, 0.8500000000000001)
Set-RowValues 10 @(8, 6799825, "Israel Premier League", "Israel Premier League", 45171.58333333334, "Hapoel Haifa", "Maccabi Petach Tikva", 2, 2, "D", 1.8, 3.25, 4, 1.95, 3.2, 3.4, -0.5, 2.1, 1.775, 2.5, 2.025, 1.825, -1, 2.2, -1, -1, 0.7749999999999999, 1.025, -1)
Set-RowValues 11 @(9, 6799829, "Israel Premier League", "Israel Premier League", 45171.58333333334, "Hapoel TelAviv", "Maccabi Netanya", 2, 0, "H", 2.4, 3.4, 2.6, 2.625, 3.4, 2.375, 0, 2, 1.85, 2.5, 2, 1.85, 1.625, -1, -1, 1, -1, -1, 0.8500000000000001)
Set-RowValues 16 @(14, 6799836, "Israel Premier League", "Israel Premier League", 45186.57291666666, "Maccabi Tel Aviv", "Maccabi Bnei Raina", 1, 1, "D", 1.181, 6, 11, 1.2, 6, 10, -1.75, 1.825, 2.025, 3, 1.85, 2, -1, 5, -1, -1, 1.025, -1, 1)
Set-RowValues 17 @(15, 6799838, "Israel Premier League", "Israel Premier League", 45186.57291666666, "Hapoel Haifa", "MS Ashdod", 2, 0, "H", 2.15, 3, 3.2, 2.15, 3.1, 3.1, -0.25, 2, 1.85, 2.5, 2, 1.85, 1.15, -1, -1, 1, -1, -1, 0.8500000000000001)
Set-RowValues 58 @(56, 7542748, "Israel Premier League", "Israel Premier League", 45276.54166666666, "MS Ashdod", "Hapoel Jerusalem FC", 2, 0, "H", 2.5, 3.2, 2.625, 2.4, 2.9, 3, -0.25, 2.125, 1.75, 2, 2.05, 1.8, 1.4, -1, -1, 1.125, -1, 0, -0)
Set-RowValues 60 @(58, 7542499, "Israel Premier League", "Israel Premier League", 45276.54166666666, "Maccabi Petach Tikva", "Hapoel Beer Sheva", 1, 4, "A", 2.65, 3.2, 2.4, 3.2, 3.3, 2.05, 0.25, 2, 1.85, 2.25, 1.85, 2, -1, -1, 1.05, -1, 0.8500000000000001, 0.8500000000000001, -1)
Set-RowValues 73 @(71, 7542719, "Israel Premier League", "Israel Premier League", 45283.54166666666, "Hapoel Haifa", "Maccabi Netanya", 2, 1, "H", 2.6, 3.1, 2.6, 2.9, 3.2, 2.3, 0.25, 1.8, 2.05, 2.5, 2, 1.85, 1.9, -1, -1, 0.8, -1, 1, -1)
Set-RowValues 74 @(72, 7542640, "Israel Premier League", "Israel Premier League", 45283.54166666666, "MS Ashdod", "Hapoel Bnei Sakhnin", 0, 1, "A", 2.05, 3.2, 3.5, 2.15, 3.1, 3.2, -0.25, 1.925, 1.925, 2.25, 1.9, 1.95, -1, -1, 2.2, -1, 0.925, -1, 0.95)
Set-RowValues 109 @(107, 7542735, "Israel Premier League", "Israel Premier League", 45304.54166666666, "Hapoel Petah Tikva", "Maccabi Netanya", 2, 0, "H", 3.75, 3.6, 1.909, 3.8, 3.75, 1.85, 0.5, 1.95, 1.9, 2.5, 1.975, 1.875, 2.8, -1, -1, 0.95, -1, -1, 0.875)
Set-RowValues 110 @(108, 7542737, "Israel Premier League", "Israel Premier League", 45304.54166666666, "MS Ashdod", "Hapoel Haifa", 0, 1, "A", 3, 3.2, 2.45, 3.2, 3.25, 2.3, 0.25, 1.85, 2, 2.25, 1.875, 1.975, -1, -1, 1.3, -1, 1, -1, 0.9750000000000001)
Set-RowValues 164 @(162, 6799984, "Israel Premier League", "Israel Premier League", 45353.5625, "Hapoel Bnei Sakhnin", "MS Ashdod", 1, 1, "D", 2, 3.25, 3.25, 2, 3.25, 3.25, -0.25, 1.825, 2.025, 2.25, 1.925, 1.925, -1, 2.25, -1, -0.5, 0.5125, -0.5, 0.4625)
Set-RowValues 165 @(163, 6799986, "Israel Premier League", "Israel Premier League", 45353.5625, "Hapoel Jerusalem FC", "Maccabi Bnei Raina", 1, 0, "H", 2.2, 3.2, 2.9, 2.3, 3.1, 2.8, -0.25, 2.1, 1.775, 2, 2.025, 1.825, 1.3, -1, -1, 1.1, -1, -1, 0.825)

# --- Fixtures that were pending (no score yet) and have now been played ---
Set-RowValues 182 @(180, 7951997, "Israel Premier League", "Israel Premier League", 45368.625, "Maccabi Tel Aviv", "Hapoel Haifa", 3, 1, "H", 1.333, 4.333, 7.5, 1.333, 4.2, 8, -1.5, 2.025, 1.825, 2.75, 1.9, 1.95, 0.333, -1, -1, 1.025, -1, 0.8999999999999999, -1)
Set-RowValues 183 @(181, 7952000, "Israel Premier League", "Israel Premier League", 45368.64583333334, "Maccabi Haifa", "Hapoel TelAviv", 0, 0, "D", 1.363, 4.333, 6.5, 1.285, 4.5, 9, -1.5, 1.95, 1.9, 2.75, 1.85, 2, -1, 3.5, -1, -1, 0.8999999999999999, -1, 1)

# --- Brand-new upcoming fixtures appended to the bottom of the table ---
Set-RowValues 184 @(182, 7986402, "Israel Premier League", "Israel Premier League", 45381.41666666666, "Hapoel Beer Sheva", "Hapoel Haifa", $null, $null, $null, 1.727, 3.4, 4.2, 1.666, 3.5, 4.5, -0.75, 1.975, 1.875, 2.25, 1.9, 1.95, 0, 0, 0, 0, 0, $null, $null)
Set-RowValues 185 @(183, 7986406, "Israel Premier League", "Israel Premier League", 45381.45833333334, "Hapoel Hadera", "Hapoel Petah Tikva", $null, $null, $null, 2.5, 3, 2.6, 2.4, 3, 2.7, 0, 1.8, 2.05, 2.25, 1.9, 1.95, 0, 0, 0, 0, 0, $null, $null)
Set-RowValues 186 @(184, 7986407, "Israel Premier League", "Israel Premier League", 45381.5625, "Maccabi Petach Tikva", "MS Ashdod", $null, $null, $null, 2.1, 3.1, 3.2, 2.15, 3.1, 3.1, -0.25, 1.925, 1.925, 2.25, 1.9, 1.95, 0, 0, 0, 0, 0, $null, $null)
Set-RowValues 187 @(185, 7986206, "Israel Premier League", "Israel Premier League", 45381.58333333334, "Maccabi Tel Aviv", "Hapoel Bnei Sakhnin", $null, $null, $null, 1.25, 6, 8, 1.285, 6, 7, -1.5, 1.85, 2, 2.75, 1.825, 2.025, 0, 0, 0, 0, 0, $null, $null)
Set-RowValues 188 @(186, 7986403, "Israel Premier League", "Israel Premier League", 45382.60416666666, "Maccabi Haifa", "Maccabi Bnei Raina", $null, $null, $null, 1.333, 4.5, 8, 1.363, 4.333, 7, -1.25, 1.925, 1.925, 2.5, 1.925, 1.925, 0, 0, 0, 0, 0, $null, $null)
